$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 22356.428
$ws.Range("J2").Value = 25249.166
$ws.Range("L2").Value = 25249.166
$ws.Range("N2").Value = -25475.166

$ws.Range("H4").Value = 100.5
$ws.Range("I4").Value = 99
$ws.Range("J4").Value = 102
$ws.Range("K4").Value = 99
$ws.Range("L4").Value = 102
$ws.Range("M4").Value = 15
$ws.Range("N4").Value = -330

$ws.Range("H5").Value = 938
$ws.Range("I5").Value = 375
$ws.Range("K5").Value = 375
$ws.Range("M5").Value = -260

$ws.Range("H13").Value = 8993
$ws.Range("J13").Value = 8993
$ws.Range("L13").Value = 8993
$ws.Range("N13").Value = -9331

$ws.Range("H17").Value = 625.087
$ws.Range("J17").Value = 625.087
$ws.Range("L17").Value = 1875.261
$ws.Range("N17").Value = -2211.261

$ws.Range("H38").Value = 4919.5
$ws.Range("J38").Value = 6492.6665
$ws.Range("L38").Value = 19477.9995
$ws.Range("N38").Value = -20221.9995

$ws.Range("H125").Value = 3966
$ws.Range("J125").Value = 2499.5
$ws.Range("L125").Value = 22495.5
$ws.Range("N125").Value = -27415.5

$ws.Range("H138").Value = 7594.3784
$ws.Range("J138").Value = 7882.2354
$ws.Range("L138").Value = 23646.7062
$ws.Range("N138").Value = -33926.7062

$ws.Range("H141").Value = 5841.222
$ws.Range("I141").Value = 5071.375
$ws.Range("K141").Value = 15214.125
$ws.Range("M141").Value = -10034.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 770.4
$ws.Range("I2").Value = 770.4
$ws.Range("K2").Value = 770.4
$ws.Range("M2").Value = -657.4

$ws.Range("H32").Value = 3058.3215
$ws.Range("I32").Value = 3023.4443
$ws.Range("K32").Value = 3023.4443
$ws.Range("M32").Value = -2736.4443

$ws.Range("H45").Value = 2863.5715
$ws.Range("I45").Value = 2388.875
$ws.Range("J45").Value = 3496.5
$ws.Range("K45").Value = 2388.875
$ws.Range("L45").Value = 3496.5
$ws.Range("M45").Value = -2011.875
$ws.Range("N45").Value = -4250.5

$ws.Range("H74").Value = 2126.389
$ws.Range("I74").Value = 636.61536
$ws.Range("K74").Value = 636.61536
$ws.Range("M74").Value = 237.38464

$ws.Range("H77").Value = 2126.389
$ws.Range("I77").Value = 636.61536
$ws.Range("K77").Value = 3183.0768
$ws.Range("M77").Value = 1184.9232

$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490

$ws.Range("H116").Value = 770.4
$ws.Range("I116").Value = 770.4
$ws.Range("K116").Value = 770.4
$ws.Range("M116").Value = 1523.6

$ws.Range("H122").Value = 1756.3684
$ws.Range("I122").Value = 1756.3684
$ws.Range("K122").Value = 5269.1052
$ws.Range("M122").Value = -2819.1052

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 770.4
$ws.Range("I3").Value = 770.4
$ws.Range("K3").Value = 770.4
$ws.Range("M3").Value = -656.4

$ws.Range("H80").Value = 1273.4
$ws.Range("J80").Value = 491.75
$ws.Range("L80").Value = 491.75
$ws.Range("N80").Value = -2487.75

$ws.Range("H83").Value = 1273.4
$ws.Range("J83").Value = 491.75
$ws.Range("L83").Value = 2458.75
$ws.Range("N83").Value = -12442.75

$ws.Range("H94").Value = 1019.8571
$ws.Range("I94").Value = 734.75
$ws.Range("J94").Value = 1400
$ws.Range("K94").Value = 734.75
$ws.Range("L94").Value = 1400
$ws.Range("M94").Value = -283.75
$ws.Range("N94").Value = -2302

$ws.Range("H105").Value = 3000
$ws.Range("I105").Value = 3000
$ws.Range("K105").Value = 3000
$ws.Range("M105").Value = -1253

$ws.Range("H125").Value = 88888
$ws.Range("J125").Value = 88888
$ws.Range("L125").Value = 88888
$ws.Range("N125").Value = -98728

$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").Value = ""

$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = ""

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""

$ws.Range("H134").Value = 3915.818
$ws.Range("I134").Value = 3687.111
$ws.Range("K134").Value = 11061.333
$ws.Range("M134").Value = -8526.332999999999

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = ""

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 175.5
$ws.Range("I7").Value = 10.5
$ws.Range("K7").Value = 10.5
$ws.Range("M7").Value = 102.5

$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").Value = ""

$ws.Range("H31").Value = 6027.85
$ws.Range("J31").Value = 9371
$ws.Range("L31").Value = 9371
$ws.Range("N31").Value = -9961

$ws.Range("H34").Value = 6027.85
$ws.Range("J34").Value = 9371
$ws.Range("L34").Value = 9371
$ws.Range("N34").Value = -9775

$ws.Range("H105").Value = 1799.6
$ws.Range("I105").Value = 1624.75
$ws.Range("J105").Value = 2499
$ws.Range("K105").Value = 1624.75
$ws.Range("L105").Value = 2499
$ws.Range("M105").Value = 122.25
$ws.Range("N105").Value = -5993

$ws.Range("H134").Value = 2605.0715
$ws.Range("I134").Value = 2372.6667
$ws.Range("J134").Value = 3999.5
$ws.Range("K134").Value = 7118.000100000001
$ws.Range("L134").Value = 11998.5
$ws.Range("M134").Value = -4583.000100000001
$ws.Range("N134").Value = -17068.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 63.25
$ws.Range("J2").Value = 80
$ws.Range("L2").Value = 480
$ws.Range("N2").Value = -706

$ws.Range("H5").Value = 887.0625
$ws.Range("J5").Value = 988.55554
$ws.Range("L5").Value = 2965.66662
$ws.Range("N5").Value = -3189.66662

$ws.Range("H17").Value = 1000
$ws.Range("I17").Value = 1000
$ws.Range("K17").Value = 3000
$ws.Range("M17").Value = -2831

$ws.Range("H135").Value = 887.0625
$ws.Range("J135").Value = 988.55554
$ws.Range("L135").Value = 8896.99986
$ws.Range("N135").Value = -13966.99986

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1068.2667
$ws.Range("I97").Value = 575.2
$ws.Range("J97").Value = 2054.4
$ws.Range("K97").Value = 575.2
$ws.Range("L97").Value = 2054.4
$ws.Range("M97").Value = -79.20000000000005
$ws.Range("N97").Value = -3046.4

$ws.Range("H122").Value = 1404.1666
$ws.Range("I122").Value = 1404.1666
$ws.Range("K122").Value = 4212.4998
$ws.Range("M122").Value = -1762.4998

$ws.Range("H126").Value = 2374.25
$ws.Range("I126").Value = 2374.25
$ws.Range("K126").Value = 7122.75
$ws.Range("M126").Value = -4652.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 4000
$ws.Range("I3").Value = 4000
$ws.Range("K3").Value = 4000
$ws.Range("M3").Value = -3888

$ws.Range("H15").Value = 4000
$ws.Range("I15").Value = 4000
$ws.Range("K15").Value = 4000
$ws.Range("M15").Value = -3830

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = ""

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").Value = ""

$ws.Range("H31").Value = 1338.3334
$ws.Range("I31").Value = 1338.3334
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1338.3334
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1090.3334
$ws.Range("N31").Value = ""

$ws.Range("H132").Value = 3020.6316
$ws.Range("I132").Value = 2278.6428
$ws.Range("K132").Value = 6835.928400000001
$ws.Range("M132").Value = -4305.928400000001

$ws.Range("H136").Value = 25955.262
$ws.Range("I136").Value = 5509.7
$ws.Range("K136").Value = 16529.1
$ws.Range("M136").Value = -13979.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1394.7
$ws.Range("I122").Value = 1305.5883
$ws.Range("K122").Value = 3916.7649
$ws.Range("M122").Value = -1466.7649

$ws.Range("H132").Value = 2409.1765
$ws.Range("I132").Value = 2120.2693
$ws.Range("K132").Value = 6360.8079
$ws.Range("M132").Value = -3830.8079

$ws.Range("H136").Value = 6008.4
$ws.Range("I136").Value = 7360.3687
$ws.Range("K136").Value = 22081.1061
$ws.Range("M136").Value = -19531.1061
